$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F59").Value = 18
$ws.Range("G59").Value = 1478.52
$ws.Range("F68").Value = 46
$ws.Range("G68").Value = 5295.52
$ws.Range("F70").Value = 13
$ws.Range("G70").Value = 1754.35
$ws.Range("F77").Value = 251
$ws.Range("G77").Value = 11731.74
$ws.Range("F85").Value = 140
$ws.Range("G85").Value = 18867.8
$ws.Range("F86").Value = 62
$ws.Range("G86").Value = 7779.14
$ws.Range("B90").Value = 178514.21
$ws.Range("F115").Value = 201
$ws.Range("G115").Value = 19458.81
$ws.Range("B117").Value = 13316.69
$ws.Range("F127").Value = 1
$ws.Range("G127").Value = 120.69
$ws.Range("B133").Value = 1844.5
$ws.Range("F164").Value = 62
$ws.Range("G164").Value = 7107.06
$ws.Range("F167").Value = 13
$ws.Range("G167").Value = 3731.13
$ws.Range("F169").Value = 5
$ws.Range("G169").Value = 717.5
$ws.Range("B175").Value = 27426.35
$ws.Range("F222").Value = 13
$ws.Range("G222").Value = 1884.09
$ws.Range("F223").Value = 15
$ws.Range("G223").Value = 1987.2
$ws.Range("F229").Value = 59
$ws.Range("G229").Value = 8465.32
$ws.Range("B243").Value = 60325
$ws.Range("E243").Value = 151.57
$ws.Range("F243").Value = -102
$ws.Range("G243").Value = -12939.72
$ws.Range("B244").Value = 63560
$ws.Range("E244").Value = 134.87
$ws.Range("F244").Value = 1
$ws.Range("G244").Value = 126.86
$ws.Range("F247").Value = 139
$ws.Range("G247").Value = 14443.49
$ws.Range("F252").Value = 0
$ws.Range("G252").Value = 0
$ws.Range("F256").Value = 282
$ws.Range("G256").Value = 42629.94
$ws.Range("B260").Value = 186140.09
$ws.Range("F270").Value = 15
$ws.Range("G270").Value = 483.6
$ws.Range("B275").Value = 5192.75
$ws.Range("F278").Value = 11
$ws.Range("G278").Value = 1510.52
$ws.Range("F282").Value = 0
$ws.Range("G282").Value = 0
$ws.Range("F283").Value = 39
$ws.Range("G283").Value = 13317.33
$ws.Range("F291").Value = 111
$ws.Range("G291").Value = 4774.11
$ws.Range("F293").Value = 33
$ws.Range("G293").Value = 2320.56
$ws.Range("F294").Value = 29
$ws.Range("G294").Value = 2069.44
$ws.Range("F299").Value = 272
$ws.Range("G299").Value = 39342.08
$ws.Range("B304").Value = 172894.39
$ws.Range("F329").Value = 27
$ws.Range("G329").Value = 4494.69
$ws.Range("B330").Value = 27397.71
$ws.Range("B380").Value = 64925
$ws.Range("E380").Value = 13.97
$ws.Range("F380").Value = 111
$ws.Range("G380").Value = 1459.65
$ws.Range("B381").Value = 45709
$ws.Range("E381").Value = 15.69
$ws.Range("F381").Value = -300
$ws.Range("G381").Value = -3945
$ws.Range("B382").Value = 64919
$ws.Range("E382").Value = 27.97
$ws.Range("F382").Value = 61
$ws.Range("G382").Value = 1604.3
$ws.Range("B383").Value = 45702
$ws.Range("E383").Value = 31.43
$ws.Range("F383").Value = -215
$ws.Range("G383").Value = -5654.5
$ws.Range("B385").Value = 53595
$ws.Range("E385").Value = 17.61
$ws.Range("F385").Value = -335
$ws.Range("G385").Value = -4934.55
$ws.Range("B386").Value = 65067
$ws.Range("E386").Value = 15.65
$ws.Range("F386").Value = 126
$ws.Range("G386").Value = 1855.98
$ws.Range("B442").Value = 64810
$ws.Range("E442").Value = 291.22
$ws.Range("F442").Value = 4
$ws.Range("G442").Value = 1095.68
$ws.Range("B443").Value = 53319
$ws.Range("E443").Value = 310.64
$ws.Range("F443").Value = -6
$ws.Range("G443").Value = -1643.52
$ws.Range("F450").Value = 10
$ws.Range("G450").Value = 1387.4
$ws.Range("F455").Value = 46
$ws.Range("G455").Value = 2926.06
$ws.Range("B460").Value = 13407.07
$ws.Range("B463").Value = 60025
$ws.Range("E463").Value = 37.22
$ws.Range("F463").Value = -98
$ws.Range("G463").Value = -3217.34
$ws.Range("B464").Value = 64833
$ws.Range("E464").Value = 34.9
$ws.Range("F464").Value = 95
$ws.Range("G464").Value = 3118.85
$ws.Range("B473").Value = 60022
$ws.Range("E473").Value = 37.22
$ws.Range("F473").Value = -113
$ws.Range("G473").Value = -3709.79
$ws.Range("B474").Value = 64830
$ws.Range("E474").Value = 34.9
$ws.Range("F474").Value = 107
$ws.Range("G474").Value = 3512.81
$ws.Range("F491").Value = 20
$ws.Range("G491").Value = 3560.4
$ws.Range("B493").Value = 11785.05
$ws.Range("F508").Value = 55
$ws.Range("G508").Value = 5716.7
$ws.Range("F509").Value = 214
$ws.Range("G509").Value = 17201.32
$ws.Range("B510").Value = 22918.02
$ws.Range("F549").Value = 25
$ws.Range("G549").Value = 1196.5
$ws.Range("F555").Value = 17
$ws.Range("G555").Value = 1182.52
$ws.Range("B560").Value = 4507.4
$ws.Range("B572").Value = 65362
$ws.Range("F572").Value = 20
$ws.Range("G572").Value = 817.4
$ws.Range("B573").Value = 65079
$ws.Range("F573").Value = 6
$ws.Range("G573").Value = 245.22
$ws.Range("F577").Value = 57
$ws.Range("G577").Value = 2450.43
$ws.Range("F578").Value = 77
$ws.Range("G578").Value = 3841.53
$ws.Range("F582").Value = 31
$ws.Range("G582").Value = 1766.69
$ws.Range("B583").Value = 16033.28
$ws.Range("F599").Value = 1612
$ws.Range("G599").Value = 262933.32
$ws.Range("F601").Value = 404
$ws.Range("G601").Value = 114279.48
$ws.Range("F602").Value = 329
$ws.Range("G602").Value = 47589.85
$ws.Range("B606").Value = 425650.7
$ws.Range("F612").Value = 31
$ws.Range("G612").Value = 1270.69
$ws.Range("F613").Value = 136
$ws.Range("G613").Value = 21645.76
$ws.Range("B618").Value = 43489.71
$ws.Range("B619").Value = 1729815.71
$ws.Range("B620").Value = 1729815.71
